# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values on the zh-cn and de-de report sheets (row 5, the 5828addd... file)
# to reflect the newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-22 03:04:38"
$wsZhCn.Range("G5").Value = "2016-01-22 03:05:26"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-22 03:04:53"
$wsDeDe.Range("G5").Value = "2016-01-22 03:05:50"
